$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn = $wb.Worksheets.Item(2)
$wsDeDe = $wb.Worksheets.Item(3)

# "Latest HO Xliff Generate Date" for cb79ddec row (shared between Overview!G4 and de-de!H4)
$wsOverview.Range("G4").Value = "2016-09-04 18:50:07"
$wsDeDe.Range("H4").Value = "2016-09-04 18:50:07"

# zh-cn sheet, cb79ddec row: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-09-04 18:49:57"
$wsZhCn.Range("K4").Value = "2016-09-04 18:50:31"

# de-de sheet, cb79ddec row: Correspond Handback DateTime
$wsDeDe.Range("K4").Value = "2016-09-04 18:50:38"
